$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# PRODUCT ADD, UPDATE API DOCUMENTATION
#
# Adds 4 new rows (13-16) to the API doc sheet describing the /product
# (GET, POST) and product/{id} (PUT, DELETE) endpoints, mirroring the
# existing CATEGORIES / BRANDS sections' layout.
# ---------------------------------------------------------------------------

# New shared-string bodies (kept as variables so we can control the exact
# order they are first written in -- that order determines the shared
# string table index they land on).
$productIdEndpoint = "product/{id}"
$productEndpoint = "/product"

$postBody = @'
{
    "name": "DROID BLACK",
    "distImportDate": "2021-05-09",
    "shopImportDate": "2021-05-09",
    "categoryId": "2",
    "variants": [
        {
        "code": "vip-droid-black",
        "mrp": "4020",
        "discount": "50",
        "size": "large",
        "colour": "black"
    }]
}
'@

$postResponse = @'
{
    "success": true,
    "data": [
        {
            "createdAt": "2021-05-09T19:10:36.976+00:00",
            "createdBy": null,
            "updatedAt": "2021-05-09T19:10:36.976+00:00",
            "updatedBy": null,
            "id": 5,
            "name": "DROID GREY",
            "distImportDate": "2021-05-09",
            "shopImportDate": "2021-05-09",
            "variants": [
                {
                    "createdAt": "2021-05-09T19:10:37.042+00:00",
                    "createdBy": null,
                    "updatedAt": "2021-05-09T19:10:37.042+00:00",
                    "updatedBy": null,
                    "id": 5,
                    "code": "vip-droid-grey",
                    "mrp": 4020,
                    "discount": 20,
                    "size": "large",
                    "colour": "grey"
                }
            ]
        }
    ]
}
'@

$putBody = @'
{
    "name": "DROID BLACK",
    "distImportDate": "2021-05-09",
    "shopImportDate": "2021-05-09",
    "categoryId": "2",
    "variants": [
        {
        "id": 2,
        "code": "vip-droid-black",
        "mrp": "4020",
        "discount": "30",
        "size": "large",
        "colour": "black"
    }]
}
'@

$putResponse = @'
{
    "success": true,
    "data": [
        {
            "createdAt": "2021-05-09T09:13:13.000+00:00",
            "createdBy": null,
            "updatedAt": "2021-05-09T19:13:53.342+00:00",
            "updatedBy": null,
            "id": 3,
            "name": "DROID BLACK",
            "distImportDate": "2021-05-09",
            "shopImportDate": "2021-05-09",
            "variants": [
                {
                    "createdAt": null,
                    "createdBy": null,
                    "updatedAt": "2021-05-09T19:13:53.343+00:00",
                    "updatedBy": null,
                    "id": 2,
                    "code": "vip-droid-black",
                    "mrp": 4020,
                    "discount": 30,
                    "size": "large",
                    "colour": "black"
                }
            ]
        }
    ]
}
'@

$putStatus = "200`n404`n400"
$postStatus = "200`n404"

# ---------------------------------------------------------------------------
# Write the brand new shared strings in the exact order they must be
# minted: product/{id}, /product, postBody, postResponse, putBody,
# putResponse, putStatus, postStatus.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = $productIdEndpoint
$ws.Range("A13").Value = $productEndpoint
$ws.Range("D14").Value = $postBody
$ws.Range("E14").Value = $postResponse
$ws.Range("D15").Value = $putBody
$ws.Range("E15").Value = $putResponse
$ws.Range("F15").Value = $putStatus
$ws.Range("F14").Value = $postStatus

# ---------------------------------------------------------------------------
# Fill in the rest of the row content (reusing already-existing shared
# strings: GET/POST/PUT/DELETE and "id").
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "GET"

$ws.Range("A14").Value = $productEndpoint
$ws.Range("B14").Value = "POST"

$ws.Range("B15").Value = "PUT"
$ws.Range("C15").Value = "id"

$ws.Range("A16").Value = $productIdEndpoint
$ws.Range("B16").Value = "DELETE"
$ws.Range("C16").Value = "id"

# ---------------------------------------------------------------------------
# Formatting: D/E columns wrap + left/top align (no horizontal for E);
# F14/F15 are the new "center/top/wrap" multi-line status-code style.
# ---------------------------------------------------------------------------
"14", "15" | ForEach-Object {
    $ws.Range("D$_").HorizontalAlignment = -4131
    $ws.Range("D$_").VerticalAlignment = -4160
    $ws.Range("D$_").WrapText = $true

    $ws.Range("E$_").VerticalAlignment = -4160
    $ws.Range("E$_").WrapText = $true

    $ws.Range("F$_").HorizontalAlignment = -4108
    $ws.Range("F$_").VerticalAlignment = -4160
    $ws.Range("F$_").WrapText = $true

    $ws.Rows.Item([int]$_).RowHeight = 409.6
}

# ---------------------------------------------------------------------------
# Update the view: scroll so row 13 is at top, select F14 (mirrors the
# author having just finished typing the new POST row's status code).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F14").Select()
